$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (employee_id, employee_name, department, absence_reason,
# absence_duration, absence_date, salary)
$data = @(
    @(43772, "Juliana da Paz", "Marketing", "Doenca", 4, 45103, 9177.549999999999),
    @(67701, "Gabriel Castro", "Atendimento ao Cliente", "Doenca", 5, 45090, 3170.61),
    @(62307, "Dr. Gael Henrique Pires", "Vendas", "Problemas pessoais", 8, 45081, 7900.74),
    @(64761, "Theo Silva", "Operacoes", "Outros", 4, 45095, 4118.28),
    @(37508, "Eduardo Carvalho", "Recursos Humanos", "Outros", 8, 45099, 4844.3),
    @(30204, "Sra. Isadora Cavalcanti", "Financeiro", "Problemas pessoais", 1, 45102, 2255.08),
    @(51898, "Ana Vitória Azevedo", "TI", "Doenca", 7, 45080, 9227.059999999999),
    @(78224, "Luna Monteiro", "Juridico", "Outros", 2, 45098, 6587.7),
    @(62566, "Sr. Anthony Gabriel da Mata", "Financeiro", "Consulta medica", 8, 45078, 6024.64),
    @(7488, "Benício Sampaio", "TI", "Consulta medica", 5, 45106, 7625.58)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $row++
}
